# Update the build timestamp embedded in the workbook text from
# "January 30 2026 16.19.47 EST" to "February 02 2026 12.49.33 EST"
# across the "About" sheet and the "Boundaries and methane sources" sheet.

$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

# --- "About" sheet: A2 and A6 contain the version/citation strings ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: mines - January 30 (built on " + $newStamp + ")"

$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Tongxin Coal Mine, China, M0363, version 'mines - January 30 (built on " + $newStamp + ")'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet: S2:S10 contain the version string ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 10; $row++) {
    $cell = $wsData.Cells.Item($row, 19)  # column S = 19
    $cell.Value = "mines - January 30 (built on " + $newStamp + ")"
}
